# Add a new worksheet that will become the "SimConfig" tab, inserted before
# the current first sheet. Excel places new sheets with the next available
# sheetId (3) and as rId1 in the package, pushing the existing sheets along.
$wb = $excel.ActiveWorkbook

$simConfig = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$simConfig.Name = "SimConfig"

# Rename the two original sheets to their satellite-specific names.
$wb.Worksheets.Item("Sheet1").Name = "Sat1"
$wb.Worksheets.Item("Sheet2").Name = "Sat2"

# Populate the new SimConfig sheet with the neural-net toggle setting.
$simConfig.Range("A1").Value = "UseNeuralNet"
$simConfig.Range("A2").Value = "no"
$simConfig.Columns.Item(1).ColumnWidth = 16.6

# Update the selection on Sat2 before switching away from it so that it is
# not left as the active/selected sheet.
$sat2 = $wb.Worksheets.Item("Sat2")
$sat2.Range("F3").Select() | Out-Null

# Make SimConfig the active sheet/tab and set its selection.
$simConfig.Activate() | Out-Null
$simConfig.Range("A10").Select() | Out-Null
